# Add the new worksheet 'MW CNN 1D' as the last sheet in the workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "MW CNN 1D"

# Header row 1: section titles ('298 Packets' / '290 Packets'), bold
$ws.Range("B1").Value = "298 Packets"
$ws.Range("B1").Font.Bold = $true
$ws.Range("F1").Value = "290 Packets"
$ws.Range("F1").Font.Bold = $true

# Header row 2: column labels, bold
$ws.Range("A2").Value = "Acc"
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Value = "Loss"
$ws.Range("B2").Font.Bold = $true
$ws.Range("C2").Value = "Time"
$ws.Range("C2").Font.Bold = $true
$ws.Range("E2").Value = "Acc"
$ws.Range("E2").Font.Bold = $true
$ws.Range("F2").Value = "Loss"
$ws.Range("F2").Font.Bold = $true
$ws.Range("G2").Value = " Time"
$ws.Range("G2").Font.Bold = $true

# Data rows 3-51 for the '298 Packets' run (Acc / Loss / Time) in columns A:C
$colA = @(89.899682998657198, 89.669644832610999, 87.739229202270494, 87.928950786590505, 87.656223773956299, 89.000874757766695, 89.169257879257202, 89.216685295104895, 89.356607198715196, 88.744753599166799, 88.946330547332707, 88.317877054214406, 89.446723461151095, 88.227760791778493, 89.000874757766695, 89.211940765380803, 88.251477479934593, 88.075983524322496, 88.993763923645005, 87.331324815750094, 87.447530031204195, 88.287049531936603, 89.098107814788804, 88.135272264480506, 88.334476947784395, 89.361351728439303, 89.306801557540894, 88.315504789352403, 87.895745038986206, 88.747125864028902, 89.370834827423096, 88.737636804580603, 88.389027118682804, 88.268077373504596, 88.8846755027771, 87.741601467132497, 89.309173822402897, 87.722629308700505, 88.891786336898804, 88.922619819641099, 88.509970903396606, 89.057791233062702, 82.941639423370304, 89.072024822235093, 88.782697916030799, 88.906013965606604, 88.336849212646399, 88.206416368484497, 88.097327947616506)
$colB = @(0.25593475341601801, 0.27932078153680101, 0.27616908617080299, 0.30643616261680501, 0.32576677804585602, 0.25954075623456002, 0.25807979315569501, 0.25581393391264901, 0.25247172338670998, 0.24845205380204, 0.27470413344004901, 0.29265578793354702, 0.25481620887200901, 0.27021681278103599, 0.25701161256516702, 0.253373333388688, 0.27868630412323597, 0.30781715617771499, 0.25344664179127702, 0.30415420046414798, 0.30682559471366999, 0.27912531950338199, 0.25731988630920799, 0.29874858677112998, 0.28577206744850697, 0.24558426476231299, 0.25280046165144898, 0.27090228943310002, 0.28100848707252302, 0.27410961023838198, 0.26262793515588601, 0.27065842287150599, 0.28455346995084202, 0.28403646209373001, 0.259713368020911, 0.286458535266501, 0.25667114041414002, 0.32560889764015499, 0.28413221820336199, 0.26338667526728399, 0.26999710821183998, 0.25949050475079799, 2.6234215325207999, 0.258677962333, 0.27610162069314398, 0.26900770298879201, 0.28707863962150998, 0.27207447418147801, 0.27564597893096798)
$colC = @(629.61373281478802, 620.31075048446598, 622.61563682556096, 620.07716345787003, 623.313329219818, 620.18612217903103, 623.250921487808, 622.42791533470097, 621.18253993988003, 623.61626076698303, 619.62087702751103, 621.20873451232899, 622.78855776786804, 620.91040968894902, 619.70832848548798, 622.00117516517605, 622.52004361152603, 620.82996916770901, 620.46274471282902, 625.07487797737099, 620.61193203925995, 621.68988919258095, 620.42281103134098, 625.19871401786804, 621.87540721893299, 624.073350429534, 625.62148094177201, 622.66134953498795, 622.57876133918705, 622.42392969131402, 622.95025801658596, 621.00262737274102, 619.62267351150501, 622.44991159438996, 625.08124709129299, 623.22033381462097, 619.80937767028797, 624.04310846328701, 622.62831640243496, 623.05654406547501, 622.40083050727799, 621.44749498367298, 620.072613477706, 624.22593569755497, 621.95561194419804, 620.71887111663796, 621.87014770507801, 619.91724228858902, 624.51877474784806)
for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 3
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
}

# Match the recorded selection on the new sheet
$ws.Range("G2").Select()
